# Update NATMI ligand-receptor edge statistics with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.09838033333333333
$ws.Range("H2").Value = 0.295141
$ws.Range("I2").Value = 0.0257774858695505
$ws.Range("J2").Value = 0.0257774858695505
$ws.Range("M2").Value = 0.6808546666666667
$ws.Range("N2").Value = 2.042564
$ws.Range("O2").Value = 0.00133071464714358
$ws.Range("P2").Value = 0.00133071464714358
$ws.Range("Q2").Value = 0.06698270905822223
$ws.Range("R2").Value = 0.602844381524
$ws.Range("S2").Value = 0.00003430247801314753
$ws.Range("T2").Value = 0.00003430247801314752
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09838033333333333
$ws.Range("H3").Value = 0.295141
$ws.Range("I3").Value = 0.0257774858695505
$ws.Range("J3").Value = 0.0257774858695505
$ws.Range("M3").Value = 310.7256466666667
$ws.Range("N3").Value = 932.1769400000001
$ws.Range("O3").Value = 0.6073060661930214
$ws.Range("P3").Value = 0.6073060661930214
$ws.Range("Q3").Value = 30.56929269428223
$ws.Range("R3").Value = 275.12363424854
$ws.Range("S3").Value = 0.01565482353978291
$ws.Range("T3").Value = 0.01565482353978291
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.09838033333333333
$ws.Range("H4").Value = 0.295141
$ws.Range("I4").Value = 0.0257774858695505
$ws.Range("J4").Value = 0.0257774858695505
$ws.Range("M4").Value = 0.149039
$ws.Range("N4").Value = 0.447117
$ws.Range("O4").Value = 0.0002912932671323377
$ws.Range("P4").Value = 0.0002912932671323377
$ws.Range("Q4").Value = 0.01466250649966667
$ws.Range("R4").Value = 0.131962558497
$ws.Range("S4").Value = 0.000007508808077399036
$ws.Range("T4").Value = 0.000007508808077399033
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.09838033333333333
$ws.Range("H5").Value = 0.295141
$ws.Range("I5").Value = 0.0257774858695505
$ws.Range("J5").Value = 0.0257774858695505
$ws.Range("M5").Value = 200.0903396666667
$ws.Range("N5").Value = 600.271019
$ws.Range("O5").Value = 0.3910719258927027
$ws.Range("P5").Value = 0.3910719258927027
$ws.Range("Q5").Value = 19.68495431318656
$ws.Range("R5").Value = 177.164588818679
$ws.Range("S5").Value = 0.01008085104367704
$ws.Range("T5").Value = 0.01008085104367704
$ws.Range("I6").Value = 0.4417479616037814
$ws.Range("J6").Value = 0.4417479616037814
$ws.Range("M6").Value = 0.6808546666666667
$ws.Range("N6").Value = 2.042564
$ws.Range("O6").Value = 0.00133071464714358
$ws.Range("P6").Value = 0.00133071464714358
$ws.Range("Q6").Value = 1.147880570623111
$ws.Range("R6").Value = 10.330925135608
$ws.Range("S6").Value = 0.0005878404828519719
$ws.Range("T6").Value = 0.0005878404828519718
$ws.Range("I7").Value = 0.4417479616037814
$ws.Range("J7").Value = 0.4417479616037814
$ws.Range("M7").Value = 310.7256466666667
$ws.Range("N7").Value = 932.1769400000001
$ws.Range("O7").Value = 0.6073060661930214
$ws.Range("P7").Value = 0.6073060661930214
$ws.Range("Q7").Value = 523.8650038916311
$ws.Range("R7").Value = 4714.78503502468
$ws.Range("S7").Value = 0.2682762168103784
$ws.Range("T7").Value = 0.2682762168103784
$ws.Range("I8").Value = 0.4417479616037814
$ws.Range("J8").Value = 0.4417479616037814
$ws.Range("M8").Value = 0.149039
$ws.Range("N8").Value = 0.447117
$ws.Range("O8").Value = 0.0002912932671323377
$ws.Range("P8").Value = 0.0002912932671323377
$ws.Range("Q8").Value = 0.2512709110193333
$ws.Range("R8").Value = 2.261438199174
$ws.Range("S8").Value = 0.000128678206984616
$ws.Range("T8").Value = 0.0001286782069846159
$ws.Range("I9").Value = 0.4417479616037814
$ws.Range("J9").Value = 0.4417479616037814
$ws.Range("M9").Value = 200.0903396666667
$ws.Range("N9").Value = 600.271019
$ws.Range("O9").Value = 0.3910719258927027
$ws.Range("P9").Value = 0.3910719258927027
$ws.Range("Q9").Value = 337.3404406511797
$ws.Range("R9").Value = 3036.063965860618
$ws.Range("S9").Value = 0.1727552261035665
$ws.Range("T9").Value = 0.1727552261035665
$ws.Range("G10").Value = 0.243138
$ws.Range("H10").Value = 0.729414
$ws.Range("I10").Value = 0.06370669977418356
$ws.Range("J10").Value = 0.06370669977418356
$ws.Range("M10").Value = 0.6808546666666667
$ws.Range("N10").Value = 2.042564
$ws.Range("O10").Value = 0.00133071464714358
$ws.Range("P10").Value = 0.00133071464714358
$ws.Range("Q10").Value = 0.165541641944
$ws.Range("R10").Value = 1.489874777496
$ws.Range("S10").Value = 0.00008477543851068469
$ws.Range("T10").Value = 0.00008477543851068467
$ws.Range("G11").Value = 0.243138
$ws.Range("H11").Value = 0.729414
$ws.Range("I11").Value = 0.06370669977418356
$ws.Range("J11").Value = 0.06370669977418356
$ws.Range("M11").Value = 310.7256466666667
$ws.Range("N11").Value = 932.1769400000001
$ws.Range("O11").Value = 0.6073060661930214
$ws.Range("P11").Value = 0.6073060661930214
$ws.Range("Q11").Value = 75.54921227924001
$ws.Range("R11").Value = 679.9429105131601
$ws.Range("S11").Value = 0.03868946522999926
$ws.Range("T11").Value = 0.03868946522999926
$ws.Range("G12").Value = 0.243138
$ws.Range("H12").Value = 0.729414
$ws.Range("I12").Value = 0.06370669977418356
$ws.Range("J12").Value = 0.06370669977418356
$ws.Range("M12").Value = 0.149039
$ws.Range("N12").Value = 0.447117
$ws.Range("O12").Value = 0.0002912932671323377
$ws.Range("P12").Value = 0.0002912932671323377
$ws.Range("Q12").Value = 0.036237044382
$ws.Range("R12").Value = 0.326133399438
$ws.Range("S12").Value = 0.00001855733271544089
$ws.Range("T12").Value = 0.00001855733271544089
$ws.Range("G13").Value = 0.243138
$ws.Range("H13").Value = 0.729414
$ws.Range("I13").Value = 0.06370669977418356
$ws.Range("J13").Value = 0.06370669977418356
$ws.Range("M13").Value = 200.0903396666667
$ws.Range("N13").Value = 600.271019
$ws.Range("O13").Value = 0.3910719258927027
$ws.Range("P13").Value = 0.3910719258927027
$ws.Range("Q13").Value = 48.649565005874
$ws.Range("R13").Value = 437.846085052866
$ws.Range("S13").Value = 0.02491390177295817
$ws.Range("T13").Value = 0.02491390177295817
$ws.Range("G14").Value = 0.6601003333333334
$ws.Range("H14").Value = 1.980301
$ws.Range("I14").Value = 0.1729586233188772
$ws.Range("J14").Value = 0.1729586233188772
$ws.Range("M14").Value = 0.6808546666666667
$ws.Range("N14").Value = 2.042564
$ws.Range("O14").Value = 0.00133071464714358
$ws.Range("P14").Value = 0.00133071464714358
$ws.Range("Q14").Value = 0.4494323924182223
$ws.Range("R14").Value = 4.044891531764001
$ws.Range("S14").Value = 0.0002301585734002191
$ws.Range("T14").Value = 0.000230158573400219
$ws.Range("G15").Value = 0.6601003333333334
$ws.Range("H15").Value = 1.980301
$ws.Range("I15").Value = 0.1729586233188772
$ws.Range("J15").Value = 0.1729586233188772
$ws.Range("M15").Value = 310.7256466666667
$ws.Range("N15").Value = 932.1769400000001
$ws.Range("O15").Value = 0.6073060661930214
$ws.Range("P15").Value = 0.6073060661930214
$ws.Range("Q15").Value = 205.1101029398823
$ws.Range("R15").Value = 1845.99092645894
$ws.Range("S15").Value = 0.1050388211419479
$ws.Range("T15").Value = 0.1050388211419479
$ws.Range("G16").Value = 0.6601003333333334
$ws.Range("H16").Value = 1.980301
$ws.Range("I16").Value = 0.1729586233188772
$ws.Range("J16").Value = 0.1729586233188772
$ws.Range("M16").Value = 0.149039
$ws.Range("N16").Value = 0.447117
$ws.Range("O16").Value = 0.0002912932671323377
$ws.Range("P16").Value = 0.0002912932671323377
$ws.Range("Q16").Value = 0.09838069357966668
$ws.Range("R16").Value = 0.8854262422170001
$ws.Range("S16").Value = 0.00005038168246526708
$ws.Range("T16").Value = 0.00005038168246526706
$ws.Range("G17").Value = 0.6601003333333334
$ws.Range("H17").Value = 1.980301
$ws.Range("I17").Value = 0.1729586233188772
$ws.Range("J17").Value = 0.1729586233188772
$ws.Range("M17").Value = 200.0903396666667
$ws.Range("N17").Value = 600.271019
$ws.Range("O17").Value = 0.3910719258927027
$ws.Range("P17").Value = 0.3910719258927027
$ws.Range("Q17").Value = 132.0796999107466
$ws.Range("R17").Value = 1188.717299196719
$ws.Range("S17").Value = 0.06763926192106383
$ws.Range("T17").Value = 0.06763926192106383
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.128962333333333
$ws.Range("H18").Value = 3.386887
$ws.Range("I18").Value = 0.2958092294336073
$ws.Range("J18").Value = 0.2958092294336073
$ws.Range("M18").Value = 0.6808546666666667
$ws.Range("N18").Value = 2.042564
$ws.Range("O18").Value = 0.00133071464714358
$ws.Range("P18").Value = 0.00133071464714358
$ws.Range("Q18").Value = 0.768659273140889
$ws.Range("R18").Value = 6.917933458268001
$ws.Range("S18").Value = 0.0003936376743675572
$ws.Range("T18").Value = 0.0003936376743675571
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.128962333333333
$ws.Range("H19").Value = 3.386887
$ws.Range("I19").Value = 0.2958092294336073
$ws.Range("J19").Value = 0.2958092294336073
$ws.Range("M19").Value = 310.7256466666667
$ws.Range("N19").Value = 932.1769400000001
$ws.Range("O19").Value = 0.6073060661930214
$ws.Range("P19").Value = 0.6073060661930214
$ws.Range("Q19").Value = 350.7975510873089
$ws.Range("R19").Value = 3157.17795978578
$ws.Range("S19").Value = 0.179646739470913
$ws.Range("T19").Value = 0.179646739470913
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1.128962333333333
$ws.Range("H20").Value = 3.386887
$ws.Range("I20").Value = 0.2958092294336073
$ws.Range("J20").Value = 0.2958092294336073
$ws.Range("M20").Value = 0.149039
$ws.Range("N20").Value = 0.447117
$ws.Range("O20").Value = 0.0002912932671323377
$ws.Range("P20").Value = 0.0002912932671323377
$ws.Range("Q20").Value = 0.1682594171976667
$ws.Range("R20").Value = 1.514334754779
$ws.Range("S20").Value = 0.00008616723688961475
$ws.Range("T20").Value = 0.00008616723688961474
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1.128962333333333
$ws.Range("H21").Value = 3.386887
$ws.Range("I21").Value = 0.2958092294336073
$ws.Range("J21").Value = 0.2958092294336073
$ws.Range("M21").Value = 200.0903396666667
$ws.Range("N21").Value = 600.271019
$ws.Range("O21").Value = 0.3910719258927027
$ws.Range("P21").Value = 0.3910719258927027
$ws.Range("Q21").Value = 225.8944567475392
$ws.Range("R21").Value = 2033.050110727853
$ws.Range("S21").Value = 0.1156826850514372
$ws.Range("T21").Value = 0.1156826850514372

Write-Output "Updated cells with new TPM values."
